$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for columns B and C (rows 2-11) ---
$bValues = @(255999, 256000, 256001, 256002, 256003, 256004, 256005, 256006, 256007, 256008)
$cValues = @(2709999, 2710000, 2710001, 2710002, 2710003, 2710004, 2710005, 2710006, 2710007, 2710008)

for ($i = 0; $i -lt 10; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
    $ws.Cells.Item($r, 3).Value = $cValues[$i]
}

# --- Row height for rows 2-11 ---
$ws.Range("A2:A11").RowHeight = 18

# --- Number format for column A (rows 2-11): station marker style "0K+000" ---
$ws.Range("A2:A11").NumberFormat = "0""K+""000"

# --- Font styling for columns B and C (rows 2-11) ---
# Build the target font once on a scratch cell, then paste the formatting onto
# the real range in a single operation to keep the style table compact.
$scratch = $ws.Range("Z1")
$scratch.Font.Name = "Helvetica Neue"
$scratch.Font.Size = 14
$scratch.Font.Color = 4671303
$scratch.Copy()
$ws.Range("B2:C11").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

# --- Update the active selection ---
[void]$ws.Range("C15").Select()
